$wb = $excel.ActiveWorkbook

# The existing sheet (currently "Tabelle1") becomes "Timesheet".
$ts = $wb.Worksheets.Item(1)
$ts.Name = "Timesheet"

# Add a new sheet before it; Worksheets.Add() inserts before the active sheet,
# which gives us the desired order: Logging, Timesheet.
$logging = $wb.Worksheets.Add()
$logging.Name = "Logging"

# Populate the new "Logging" sheet with its carryover bookkeeping values.
$logging.Range("B1").Value = "carryover"
$logging.Range("A2").Value = "row"
$logging.Range("B2").Value = 36
$logging.Range("A3").Value = "column"
$logging.Range("B3").Value = 10

# Match the view state captured in the saved file: Logging's selection sits
# on B2 at a slightly larger zoom level.
$logging.Range("B2").Select()
$excel.ActiveWindow.Zoom = 110

# Timesheet remains the active/visible tab, with its selection moved to J36.
# Re-fetch by name: the reference grabbed before Add() can be stale afterwards.
$ts = $wb.Worksheets.Item("Timesheet")
$ts.Activate()
$ts.Range("J36").Select()
